$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 (school_district): add missing libraries_tools value
$ws.Range("H13").Value = "pandas"

# Row 14 (dev_tool_cli): add missing repo hyperlink (tools.py)
$ws.Hyperlinks.Add($ws.Range("F14"), "https://github.com/cdpeters/portfolio-website/blob/main/tools.py", "", "", "https://github.com/cdpeters/portfolio-website/blob/main/tools.py")

# Row 15 (dash_test_app): fix repo link to point at the correct repo
$ws.Range("F15").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F15"), "https://github.com/cdpeters/dash-test-app", "", "", "https://github.com/cdpeters/dash-test-app")

# Row 23 (guides): add missing repo hyperlink (notes folder)
$ws.Hyperlinks.Add($ws.Range("F23"), "https://github.com/cdpeters/portfolio-website/tree/main/notes", "", "", "https://github.com/cdpeters/portfolio-website/tree/main/notes")
